{"js": "// Merge the three \"Installation\" sentences into a single run and append\n// two new sentences about the desktop shortcut.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst newText =\n  \"Double-click the setup.exe file to install the necessary components for the program to operate correctly.  \" +\n  \"The program may run after the installation of the Microsoft SQL software bundle.  \" +\n  \"The first application window that the user will see is a login page. \" +\n  \"A shortcut will be added to the desktop to access the program.  \" +\n  \"The shortcut may be moved. \";\n\nlet installPara = null;\nlet emptyBorderPara = null;\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const p = paras.items[i];\n  if (p.text.indexOf(\"Double-click the setup.exe file\") !== -1) {\n    installPara = p;\n  }\n  if (p.text.indexOf(\"Save: saves a record to the database table.\") !== -1) {\n    // The paragraph immediately following this one is the empty bordered\n    // paragraph that must be removed.\n    emptyBorderPara = paras.items[i + 1];\n  }\n}\n\nif (installPara) {\n  installPara.insertText(newText, \"Replace\");\n}\n\nif (emptyBorderPara) {\n  emptyBorderPara.delete();\n}\n\nawait context.sync();\n", "ps1": "# Merge the three \"Installation\" sentences into a single run and append\n# two new sentences about the desktop shortcut, then remove the now\n# redundant empty (bordered) paragraph that used to sit right after the\n# \"Save: saves a record to the database table.\" line.\n\n$d = $word.ActiveDocument\n\n$newText = \"Double-click the setup.exe file to install the necessary components for the program to operate correctly.  The program may run after the installation of the Microsoft SQL software bundle.  The first application window that the user will see is a login page. A shortcut will be added to the desktop to access the program.  The shortcut may be moved. \"\n\n$installIndex = -1\n$emptyIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Double-click the setup.exe file*\") {\n        $installIndex = $i\n    }\n    if ($t -like \"*Save: saves a record to the database table.*\") {\n        $emptyIndex = $i + 1\n    }\n}\n\nif ($installIndex -gt 0) {\n    $p = $d.Paragraphs.Item($installIndex)\n    $r = $d.Range($p.Range.Start, $p.Range.End)\n    $r.Text = $newText\n}\n\nif ($emptyIndex -gt 0) {\n    $d.Paragraphs.Item($emptyIndex).Range.Delete()\n}\n"}
